# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures on Sheet1 to reflect the latest scrape, as produced by the
# GitHub Actions job that regenerates cryptos.xlsx.
#
# Both columns store plain text (not numbers) in the source data - e.g.
# "46.984.39" or "0.998" are text labels, and "  +4.31%  " keeps its
# padding spaces. Forcing the cell's NumberFormat to "@" (Text) before
# assigning the value stops Excel's automatic number/date coercion, and
# resetting the Style back to "Normal" afterwards keeps the cell's
# formatting identical to the untouched cells around it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
  $r = $ws.Range($addr)
  $r.NumberFormat = '@'
  $r.Value = $val
  $r.Style = 'Normal'
}

Set-TextValue 'D2' '46.984.39'
Set-TextValue 'E2' '  +4.31%  '
Set-TextValue 'D3' '2.491.88'
Set-TextValue 'E3' '  +2.68%  '
Set-TextValue 'D4' '0.998'
Set-TextValue 'E4' '  -0.20%  '
Set-TextValue 'D5' '322.72'
Set-TextValue 'E5' '  +1.41%  '
Set-TextValue 'D6' '105.26'
Set-TextValue 'E6' '  +0.93%  '
Set-TextValue 'E7' '  +0.79%  '
Set-TextValue 'D8' '0.998'
Set-TextValue 'E8' '  -0.28%  '
Set-TextValue 'D9' '0.539'
Set-TextValue 'E9' '  +1.62%  '
Set-TextValue 'D10' '37.01'
Set-TextValue 'E10' '  +3.61%  '
Set-TextValue 'D11' '0.0812'
Set-TextValue 'E11' '  +1.17%  '
Set-TextValue 'E12' '  +0.50%  '
Set-TextValue 'D13' '18.31'
Set-TextValue 'E13' '  -1.30%  '
Set-TextValue 'D14' '7.24'
Set-TextValue 'E14' '  +4.02%  '
Set-TextValue 'D15' '2.876.61'
Set-TextValue 'E15' '  +2.32%  '
Set-TextValue 'D16' '2.517.17'
Set-TextValue 'E16' '  +3.65%  '
Set-TextValue 'D17' '0.843'
Set-TextValue 'E17' '  +0.99%  '
Set-TextValue 'D18' '46.838.60'
Set-TextValue 'E18' '  +4.26%  '
Set-TextValue 'D19' '12.66'
Set-TextValue 'E19' '  +2.19%  '
Set-TextValue 'D20' '6.61'
Set-TextValue 'E20' '  +4.05%  '
Set-TextValue 'D21' '0.0₃0934'
Set-TextValue 'E21' '  +1.58%  '
Set-TextValue 'D22' '70.64'
Set-TextValue 'E22' '  +2.62%  '
Set-TextValue 'D23' '251.05'
Set-TextValue 'E23' '  +3.06%  '
Set-TextValue 'E24' '  +3.16%  '
Set-TextValue 'D25' '2.55'
Set-TextValue 'E25' '  +1.46%  '
Set-TextValue 'D26' '26.19'
Set-TextValue 'E26' '  +3.11%  '
Set-TextValue 'D27' '0.999'
Set-TextValue 'E27' '  -0.22%  '
Set-TextValue 'D28' '10.03'
Set-TextValue 'E28' '  +4.69%  '
Set-TextValue 'E29' '  +0.29%  '
Set-TextValue 'D30' '34.84'
Set-TextValue 'E30' '  +2.81%  '
Set-TextValue 'E31' '  +3.68%  '
Set-TextValue 'D32' '49.54'
Set-TextValue 'E32' '  +1.32%  '
Set-TextValue 'D33' '19.65'
Set-TextValue 'E33' '  -0.58%  '
Set-TextValue 'D34' '5.32'
Set-TextValue 'E34' '  +1.77%  '
Set-TextValue 'D35' '0.0778'
Set-TextValue 'E35' '  +1.74%  '
Set-TextValue 'D36' '1.00'
Set-TextValue 'E36' '  +0.02%  '
Set-TextValue 'E37' '  +1.27%  '
Set-TextValue 'D38' '4.60'
Set-TextValue 'E38' '  +1.97%  '
Set-TextValue 'D39' '2.96'
Set-TextValue 'E39' '  +2.80%  '
Set-TextValue 'D40' '122.80'
Set-TextValue 'E40' '  -3.14%  '
Set-TextValue 'E41' '  +1.47%  '
Set-TextValue 'E42' '  +2.12%  '
Set-TextValue 'D43' '21.48'
Set-TextValue 'E43' '  +1.84%  '
Set-TextValue 'E44' '  +1.53%  '
Set-TextValue 'D45' '1.961.92'
Set-TextValue 'E45' '  +0.96%  '
Set-TextValue 'D46' '2.99'
Set-TextValue 'E46' '  +1.05%  '
Set-TextValue 'E47' '  +0.58%  '
Set-TextValue 'D48' '1.79'
Set-TextValue 'E48' '  +0.46%  '
Set-TextValue 'E49' '  -1.60%  '
Set-TextValue 'E50' '  +15.45%  '
Set-TextValue 'D51' '79.05'
Set-TextValue 'E51' '  +4.29%  '
